$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 50.70817566666667
$ws.Range("H2").Value = 152.124527
$ws.Range("I2").Value = 0.5661129211027078
$ws.Range("J2").Value = 0.5661129211027077
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.441269
$ws.Range("N2").Value = 40.323807
$ws.Range("O2").Value = 0.0897308213348123
$ws.Range("P2").Value = 0.08973082133481232
$ws.Range("Q2").Value = 681.5822296349211
$ws.Range("R2").Value = 6134.240066714289
$ws.Range("S2").Value = 0.05079777737879577
$ws.Range("T2").Value = 0.05079777737879577

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 50.70817566666667
$ws.Range("H3").Value = 152.124527
$ws.Range("I3").Value = 0.5661129211027078
$ws.Range("J3").Value = 0.5661129211027077
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 54.711535
$ws.Range("N3").Value = 164.134605
$ws.Range("O3").Value = 0.3652416280068742
$ws.Range("P3").Value = 0.3652416280068742
$ws.Range("Q3").Value = 2774.322127772982
$ws.Range("R3").Value = 24968.89914995683
$ws.Range("S3").Value = 0.2067680049392801
$ws.Range("T3").Value = 0.2067680049392801

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 50.70817566666667
$ws.Range("H4").Value = 152.124527
$ws.Range("I4").Value = 0.5661129211027078
$ws.Range("J4").Value = 0.5661129211027077
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 63.67711
$ws.Range("N4").Value = 191.03133
$ws.Range("O4").Value = 0.4250937452800914
$ws.Range("P4").Value = 0.4250937452800915
$ws.Range("Q4").Value = 3228.950079825657
$ws.Range("R4").Value = 29060.55071843091
$ws.Range("S4").Value = 0.240651061883003
$ws.Range("T4").Value = 0.2406510618830029

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 50.70817566666667
$ws.Range("H5").Value = 152.124527
$ws.Range("I5").Value = 0.5661129211027078
$ws.Range("J5").Value = 0.5661129211027077
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.96553866666667
$ws.Range("N5").Value = 53.896616
$ws.Range("O5").Value = 0.119933805378222
$ws.Range("P5").Value = 0.119933805378222
$ws.Range("Q5").Value = 910.9996906556258
$ws.Range("R5").Value = 8198.997215900632
$ws.Range("S5").Value = 0.06789607690162891
$ws.Range("T5").Value = 0.0678960769016289

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.08683666666667
$ws.Range("H6").Value = 51.26051
$ws.Range("I6").Value = 0.1907597520636141
$ws.Range("J6").Value = 0.1907597520636141
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.441269
$ws.Range("N6").Value = 40.323807
$ws.Range("O6").Value = 0.0897308213348123
$ws.Range("P6").Value = 0.08973082133481232
$ws.Range("Q6").Value = 229.66876799573
$ws.Range("R6").Value = 2067.01891196157
$ws.Range("S6").Value = 0.01711702923029325
$ws.Range("T6").Value = 0.01711702923029325

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.08683666666667
$ws.Range("H7").Value = 51.26051
$ws.Range("I7").Value = 0.1907597520636141
$ws.Range("J7").Value = 0.1907597520636141
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 54.711535
$ws.Range("N7").Value = 164.134605
$ws.Range("O7").Value = 0.3652416280068742
$ws.Range("P7").Value = 0.3652416280068742
$ws.Range("Q7").Value = 934.8470623276166
$ws.Range("R7").Value = 8413.62356094855
$ws.Range("S7").Value = 0.06967340240190208
$ws.Range("T7").Value = 0.06967340240190208

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.08683666666667
$ws.Range("H8").Value = 51.26051
$ws.Range("I8").Value = 0.1907597520636141
$ws.Range("J8").Value = 0.1907597520636141
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 63.67711
$ws.Range("N8").Value = 191.03133
$ws.Range("O8").Value = 0.4250937452800914
$ws.Range("P8").Value = 0.4250937452800915
$ws.Range("Q8").Value = 1088.040377975367
$ws.Range("R8").Value = 9792.3634017783
$ws.Range("S8").Value = 0.08109077745342334
$ws.Range("T8").Value = 0.08109077745342336

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.08683666666667
$ws.Range("H9").Value = 51.26051
$ws.Range("I9").Value = 0.1907597520636141
$ws.Range("J9").Value = 0.1907597520636141
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.96553866666667
$ws.Range("N9").Value = 53.896616
$ws.Range("O9").Value = 0.119933805378222
$ws.Range("P9").Value = 0.119933805378222
$ws.Range("Q9").Value = 306.9742248260178
$ws.Range("R9").Value = 2762.76802343416
$ws.Range("S9").Value = 0.02287854297799537
$ws.Range("T9").Value = 0.02287854297799537

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 18.33915266666667
$ws.Range("H10").Value = 55.017458
$ws.Range("I10").Value = 0.2047407770084672
$ws.Range("J10").Value = 0.2047407770084672
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.441269
$ws.Range("N10").Value = 40.323807
$ws.Range("O10").Value = 0.0897308213348123
$ws.Range("P10").Value = 0.08973082133481232
$ws.Range("Q10").Value = 246.501484224734
$ws.Range("R10").Value = 2218.513358022606
$ws.Range("S10").Value = 0.01837155808169741
$ws.Range("T10").Value = 0.01837155808169742

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 18.33915266666667
$ws.Range("H11").Value = 55.017458
$ws.Range("I11").Value = 0.2047407770084672
$ws.Range("J11").Value = 0.2047407770084672
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 54.711535
$ws.Range("N11").Value = 164.134605
$ws.Range("O11").Value = 0.3652416280068742
$ws.Range("P11").Value = 0.3652416280068742
$ws.Range("Q11").Value = 1003.363192992677
$ws.Range("R11").Value = 9030.268736934091
$ws.Range("S11").Value = 0.07477985471396494
$ws.Range("T11").Value = 0.07477985471396495

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 18.33915266666667
$ws.Range("H12").Value = 55.017458
$ws.Range("I12").Value = 0.2047407770084672
$ws.Range("J12").Value = 0.2047407770084672
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 63.67711
$ws.Range("N12").Value = 191.03133
$ws.Range("O12").Value = 0.4250937452800914
$ws.Range("P12").Value = 0.4250937452800915
$ws.Range("Q12").Value = 1167.784241662127
$ws.Range("R12").Value = 10510.05817495914
$ws.Range("S12").Value = 0.08703402371008534
$ws.Range("T12").Value = 0.08703402371008534

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 18.33915266666667
$ws.Range("H13").Value = 55.017458
$ws.Range("I13").Value = 0.2047407770084672
$ws.Range("J13").Value = 0.2047407770084672
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.96553866666667
$ws.Range("N13").Value = 53.896616
$ws.Range("O13").Value = 0.119933805378222
$ws.Range("P13").Value = 0.119933805378222
$ws.Range("Q13").Value = 329.4727563469031
$ws.Range("R13").Value = 2965.254807122129
$ws.Range("S13").Value = 0.02455534050271945
$ws.Range("T13").Value = 0.02455534050271945

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.438381
$ws.Range("H14").Value = 10.315143
$ws.Range("I14").Value = 0.03838654982521095
$ws.Range("J14").Value = 0.03838654982521095
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 13.441269
$ws.Range("N14").Value = 40.323807
$ws.Range("O14").Value = 0.0897308213348123
$ws.Range("P14").Value = 0.08973082133481232
$ws.Range("Q14").Value = 46.216203945489
$ws.Range("R14").Value = 415.945835509401
$ws.Range("S14").Value = 0.003444456644025875
$ws.Range("T14").Value = 0.003444456644025875

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.438381
$ws.Range("H15").Value = 10.315143
$ws.Range("I15").Value = 0.03838654982521095
$ws.Range("J15").Value = 0.03838654982521095
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 54.711535
$ws.Range("N15").Value = 164.134605
$ws.Range("O15").Value = 0.3652416280068742
$ws.Range("P15").Value = 0.3652416280068742
$ws.Range("Q15").Value = 188.119102424835
$ws.Range("R15").Value = 1693.071921823515
$ws.Range("S15").Value = 0.01402036595172704
$ws.Range("T15").Value = 0.01402036595172704

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.438381
$ws.Range("H16").Value = 10.315143
$ws.Range("I16").Value = 0.03838654982521095
$ws.Range("J16").Value = 0.03838654982521095
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 63.67711
$ws.Range("N16").Value = 191.03133
$ws.Range("O16").Value = 0.4250937452800914
$ws.Range("P16").Value = 0.4250937452800915
$ws.Range("Q16").Value = 218.94616515891
$ws.Range("R16").Value = 1970.51548643019
$ws.Range("S16").Value = 0.01631788223357976
$ws.Range("T16").Value = 0.01631788223357977

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.438381
$ws.Range("H17").Value = 10.315143
$ws.Range("I17").Value = 0.03838654982521095
$ws.Range("J17").Value = 0.03838654982521095
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.96553866666667
$ws.Range("N17").Value = 53.896616
$ws.Range("O17").Value = 0.119933805378222
$ws.Range("P17").Value = 0.119933805378222
$ws.Range("Q17").Value = 61.77236680623201
$ws.Range("R17").Value = 555.9513012560881
$ws.Range("S17").Value = 0.004603844995878272
$ws.Range("T17").Value = 0.004603844995878273

Write-Output "applied edits"